$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LH_TC_NOTIFICATION_REVIEWS")
$ws2 = $wb.Worksheets.Item("Version History")

# --- Version History sheet: log v2.0 review + v2.1 closure ---
$ws2.Range("A5").Value = "V2.0"
$ws2.Range("B5").Value = "Mahmoud Abdelmageed"
$ws2.Range("C5").Value = "Reviewed v2.0"
$ws2.Range("D5").Value = 45789

$ws2.Range("A6").Value = "V2.1"
$ws2.Range("B6").Value = "Mahmoud Abdelmageed"
$ws2.Range("C6").Value = "closed reviewer status"
$ws2.Range("D6").Value = 45789

# --- Notification reviews sheet: new review row for v2.0 ---
$ws1.Range("C5").Value = "LH-TC-NOTIFICATION-005`nLH-TC-NOTIFICATION-003"
$ws1.Range("B5").Value = "LH-TC-NOTIGICATION-REVIEW-004"
$ws1.Range("B5").VerticalAlignment = -4108
$ws1.Range("A5").Value = 45789
$ws1.Range("D5").Value = "Mahmoud Abdelmageed"
$ws1.Range("E5").Value = "v2.0"
$ws1.Range("F5").Value = "The notifications are not a dropdown"
$ws1.Range("G5").Value = "Adjust testcase so that the notifications are a section not a dropdown"
$ws1.Range("H5").Value = "Hala Eldaly"
$ws1.Range("I5").Value = "Closed"
$ws1.Range("J5").Value = "Closed"

$ws1.Rows.Item(5).RowHeight = 42

# --- Window / selection state ---
$ws2.Activate()
$ws2.Range("B14").Select()
$ws1.Activate()
$ws1.Range("D10").Select()
